$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- UC (use case) column updates -----------------------------------------
# F3 / H3: "Disponibilizar Solução" -> "Publicar Solução de Problemas"
$ws.Range("F3").Value = "Publicar Solução de Problemas"
$ws.Range("H3").Value = "Publicar Solução de Problemas"

# H9: "Publicar Solução de Problemas" -> "Logar no Sistema"
$ws.Range("H9").Value = "Logar no Sistema"

# H18: "Logar no Sistema" -> "Consultar Problemas" (new use case diagram entry)
$ws.Range("H18").Value = "Consultar Problemas"

# F3 / H3 pick up the same vertical-centre alignment used by the rest of the
# "UC" table (style index 4 => <alignment vertical="center"/>)
$ws.Range("F3").VerticalAlignment = -4108
$ws.Range("H3").VerticalAlignment = -4108

# --- Clear the old sequential numbering in column E (rows 3-20) -----------
$ws.Range("E3:E20").Value = $null

# --- Column G (rows 12-18) gets the same vertical-centre alignment --------
$ws.Range("G12:G18").VerticalAlignment = -4108

# --- Update the active selection to G20 (was G19) --------------------------
$ws.Range("G20").Select() | Out-Null
